# ============================================================================
# experiments.xlsx - "experiments table sorted by categories"
# ============================================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ----------------------------------------------------------------------------
# 1. Workbook-level metadata
# ----------------------------------------------------------------------------
# Print area: top row moves from 42 to 24
$ws.PageSetup.PrintArea = '$D$24:$F$55'

# ----------------------------------------------------------------------------
# 2. Re-order / rewrite the category summary table (was rows 24-38, now 22-36)
#    Row 22 = header, rows 23-36 = the 14 category totals, sorted.
# ----------------------------------------------------------------------------

# -- Data rows (values only, no new shared strings introduced here) --------
$ws.Range("A23").Value = "Single"
$ws.Range("B23").Value = 8
$ws.Range("C23").Value = 10
$ws.Range("D23").Value = 7
$ws.Range("E23").Value = 16
$ws.Range("F23").Value = 11
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = 10
$ws.Range("I23").Value = 9
$ws.Range("J23").Value = 2
$ws.Range("K23").Value = 9
$ws.Range("L23").Value = 3
$ws.Range("M23").Value = 9
$ws.Range("N23").Value = 4
$ws.Range("O23").Value = 9
$ws.Range("P23").Value = 11
$ws.Range("Q23").Formula = "=SUM(B23:P23)"

$ws.Range("A24").Value = "Open-ended"
$ws.Range("B24").Value = 11
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 11
$ws.Range("G24").Value = 6
$ws.Range("H24").Value = 5
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = 6
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 12
$ws.Range("M24").Value = 3
$ws.Range("N24").Value = 6
$ws.Range("O24").Value = 9
$ws.Range("P24").Value = 2
$ws.Range("Q24").Formula = "=SUM(B24:P24)"

$ws.Range("A25").Value = "Multiple"
$ws.Range("B25").Value = 3
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 1
$ws.Range("I25").Value = 6
$ws.Range("J25").Value = 8
$ws.Range("K25").Value = 3
$ws.Range("L25").Value = 8
$ws.Range("M25").Value = 2
$ws.Range("N25").Value = 1
$ws.Range("O25").Value = 7
$ws.Range("P25").Value = 8
$ws.Range("Q25").Formula = "=SUM(B25:P25)"

$ws.Range("A26").Value = "Section"
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 7
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 2
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = 3
$ws.Range("K26").Value = 5
$ws.Range("L26").Value = 3
$ws.Range("M26").Value = 2
$ws.Range("N26").Value = 1
$ws.Range("O26").Value = 1
$ws.Range("P26").Value = 1
$ws.Range("Q26").Formula = "=SUM(B26:P26)"

$ws.Range("A27").Value = "Grid"
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = 1
$ws.Range("K27").Value = 2
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 5
$ws.Range("N27").Value = 6
$ws.Range("O27").Value = 7
$ws.Range("P27").Value = 1
$ws.Range("Q27").Formula = "=SUM(B27:P27)"

$ws.Range("A28").Value = "Intro"
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 1
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 1
$ws.Range("I28").Value = 2
$ws.Range("J28").Value = 1
$ws.Range("K28").Value = 1
$ws.Range("L28").Value = 3
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 1
$ws.Range("O28").Value = 1
$ws.Range("P28").Value = 1
$ws.Range("Q28").Formula = "=SUM(B28:P28)"

$ws.Range("A29").Value = "Filter"
$ws.Range("B29").Value = 1
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 9
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 3
$ws.Range("J29").Value = 9
$ws.Range("K29").Value = 3
$ws.Range("L29").Value = 6
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 5
$ws.Range("O29").Value = 5
$ws.Range("P29").Value = 6
$ws.Range("Q29").Formula = "=SUM(B29:P29)"

$ws.Range("A30").Value = "Skip logic"
$ws.Range("B30").Value = 3
$ws.Range("C30").Value = 4
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 4
$ws.Range("I30").Value = 2
$ws.Range("J30").Value = 2
$ws.Range("K30").Value = 4
$ws.Range("L30").Value = 2
$ws.Range("M30").Value = 4
$ws.Range("N30").Value = 1
$ws.Range("O30").Value = 1
$ws.Range("P30").Value = 5
$ws.Range("Q30").Formula = "=SUM(B30:P30)"

$ws.Range("A31").Value = "Loop"
$ws.Range("B31").Value = 0
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 2
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 6
$ws.Range("M31").Value = 0
$ws.Range("N31").Value = 0
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = 0
$ws.Range("Q31").Formula = "=SUM(B31:P31)"

$ws.Range("A32").Value = "Computation"
$ws.Range("B32").Value = 0
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("N32").Value = 4
$ws.Range("O32").Value = 0
$ws.Range("P32").Value = 0
$ws.Range("Q32").Formula = "=SUM(B32:P32)"

$ws.Range("A33").Value = "Check"
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = 0
$ws.Range("Q33").Formula = "=SUM(B33:P33)"

$ws.Range("A34").Value = "Piping"
$ws.Range("B34").Value = 0
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 3
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 4
$ws.Range("J34").Value = 9
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 1
$ws.Range("O34").Value = 0
$ws.Range("P34").Value = 5
$ws.Range("Q34").Formula = "=SUM(B34:P34)"

# -- These two rows introduce brand-new category labels ("Randomising" /
#    "Rotating") - written before the new "Feature" header label so that
#    the shared-string table grows in the same order as the source edit. --
$ws.Range("A35").Value = "Randomising"
$ws.Range("B35").Value = 0
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 2
$ws.Range("J35").Value = 1
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 0
$ws.Range("N35").Value = 0
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = 3
$ws.Range("Q35").Formula = "=SUM(B35:P35)"

$ws.Range("A36").Value = "Rotating"
$ws.Range("B36").Value = 0
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 2
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("O36").Value = 0
$ws.Range("P36").Value = 2
$ws.Range("Q36").Formula = "=SUM(B36:P36)"

# -- Header row (row 22) - moved up from its old position at row 24, and
#    now carries a label in column A ("Feature") --------------------------
$ws.Range("A22").Font.Bold = $true
$ws.Range("A22").Value = "Feature"
$ws.Range("B22:Q22").Font.Bold = $true
$ws.Range("B22:Q22").NumberFormat = "@"
$ws.Range("B22").Value = "01"
$ws.Range("C22").Value = "02"
$ws.Range("D22").Value = "03"
$ws.Range("E22").Value = "04"
$ws.Range("F22").Value = "05"
$ws.Range("G22").Value = "06"
$ws.Range("H22").Value = "07"
$ws.Range("I22").Value = "08"
$ws.Range("J22").Value = "09"
$ws.Range("K22").Value = "10"
$ws.Range("L22").Value = "11"
$ws.Range("M22").Value = "12"
$ws.Range("N22").Value = "13"
$ws.Range("O22").Value = "14"
$ws.Range("P22").Value = "15"
$ws.Range("Q22").Value = "TOTAL"

# -- Old rows 37 & 38 (previously "Single" / "Skip logic") are now blank,
#    leaving only the pre-existing "V" helper-column cell in each row. ----
$ws.Range("A37:Q38").Clear()

# ----------------------------------------------------------------------------
# 3. Drop the stray formatted placeholder rows 42-54 below the table
#    (row 55 is left untouched).
# ----------------------------------------------------------------------------
$ws.Range("F42:F54").Clear()

# ----------------------------------------------------------------------------
# 4. View state: scroll position & selection
# ----------------------------------------------------------------------------
$ws.Range("A22:Q36").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$ws.Range("Q36").Activate()
